# Update column G ("K") values on the active worksheet to reflect the
# regenerated save_data (K replaces old Strike# values; std/mean recalculated
# and s_vals rewritten upstream). Only column G values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gUpdates = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 2
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 2
    14 = 2
    15 = 0
    16 = 2
    17 = 2
    18 = 1
    21 = 1
    23 = 1
}

foreach ($row in $gUpdates.Keys) {
    $ws.Range("G$row").Value = $gUpdates[$row]
}
